$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = "02/05/2024"
$ws.Range("B14").Value = "SO240205001"
$ws.Range("C14").Value = "a"
$ws.Range("D14").Value = "7539514563"
$ws.Range("E14").Value = "a"
$ws.Range("F14").Value = "a"
$ws.Range("G14").Value = "20.0"
$ws.Range("H14").Value = "50.0"
$ws.Range("I14").Value = "AMS"
$ws.Range("J14").Value = "DVD"
$ws.Range("K14").Value = "Ashley"
$ws.Range("L14").Value = "YES"
$ws.Range("M14").Value = "1234 Address"
$ws.Range("N14").Value = "CITY"
$ws.Range("O14").Value = "TX"
$ws.Range("P14").Value = "79935"

# Row 16
$ws.Range("A16").Value = "02/14/2024"
$ws.Range("B16").Value = "SO240214001"
$ws.Range("C16").Value = "AB"
$ws.Range("D16").Value = "9999999999"
$ws.Range("E16").Value = "ARTIST"
$ws.Range("F16").Value = "TITLE"
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = "AEC"
$ws.Range("J16").Value = "LP"
$ws.Range("K16").Value = "Ashley"
$ws.Range("L16").Value = "YES"
$ws.Range("M16").Value = "1234 test address"
$ws.Range("N16").Value = "city"
$ws.Range("O16").Value = "tx"
$ws.Range("P16").Value = "79935"
